{"js": "// Update the date heading above the table.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2025-09-04 Thursday\", Word.InsertLocation.replace);\n\n// Update the division-facts table. The table is made of 5 \"data\" rows,\n// each one followed by 3 blank spacer rows, 5 cells (columns) per row.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Index (within `rows.items`) of the 5 rows that actually hold text.\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\nfor (const idx of dataRowIndexes) {\n  const cells = rows.items[idx].cells;\n  cells.load(\"items\");\n}\nawait context.sync();\n\n// New text for every cell, per data row, left-to-right.\nconst newValues = [\n  [\"94\u00f78=11, 6\", \"68\u00f75=13, 3\", \"88\u00f75=17, 3\", \"72\u00f72=36, 0\", \"22\u00f79=2, 4\"],\n  [\"66\u00f77=9, 3\", \"62\u00f79=6, 8\", \"77\u00f79=8, 5\", \"19\u00f73=6, 1\", \"88\u00f75=17, 3\"],\n  [\"64\u00f76=10, 4\", \"41\u00f77=5, 6\", \"17\u00f79=1, 8\", \"40\u00f74=10, 0\", \"57\u00f79=6, 3\"],\n  [\"17\u00f77=2, 3\", \"33\u00f78=4, 1\", \"49\u00f74=12, 1\", \"35\u00f78=4, 3\", \"82\u00f72=41, 0\"],\n  [\"31\u00f72=15, 1\", \"97\u00f74=24, 1\", \"65\u00f76=10, 5\", \"75\u00f72=37, 1\", \"37\u00f78=4, 5\"],\n];\n\ndataRowIndexes.forEach((idx, rowPos) => {\n  const cells = rows.items[idx].cells;\n  const values = newValues[rowPos];\n  for (let c = 0; c < cells.items.length; c++) {\n    cells.items[c].value = values[c];\n  }\n});\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading above the table.\n$d.Paragraphs.Item(1).Range.Text = \"2025-09-04 Thursday\"\n\n# Update the division-facts table. The table is made of 5 \"data\" rows,\n# each one followed by 3 blank spacer rows, 5 cells (columns) per row.\n$t = $d.Tables.Item(1)\n\n# 1-based row indexes of the 5 rows that actually hold text.\n$dataRowIndexes = @(1, 5, 9, 13, 17)\n\n# New text for every cell, per data row, left-to-right.\n$newValues = @(\n  @(\"94\u00f78=11, 6\", \"68\u00f75=13, 3\", \"88\u00f75=17, 3\", \"72\u00f72=36, 0\", \"22\u00f79=2, 4\"),\n  @(\"66\u00f77=9, 3\", \"62\u00f79=6, 8\", \"77\u00f79=8, 5\", \"19\u00f73=6, 1\", \"88\u00f75=17, 3\"),\n  @(\"64\u00f76=10, 4\", \"41\u00f77=5, 6\", \"17\u00f79=1, 8\", \"40\u00f74=10, 0\", \"57\u00f79=6, 3\"),\n  @(\"17\u00f77=2, 3\", \"33\u00f78=4, 1\", \"49\u00f74=12, 1\", \"35\u00f78=4, 3\", \"82\u00f72=41, 0\"),\n  @(\"31\u00f72=15, 1\", \"97\u00f74=24, 1\", \"65\u00f76=10, 5\", \"75\u00f72=37, 1\", \"37\u00f78=4, 5\")\n)\n\nfor ($i = 0; $i -lt $dataRowIndexes.Count; $i++) {\n  $row = $t.Rows.Item($dataRowIndexes[$i])\n  $values = $newValues[$i]\n  for ($c = 1; $c -le $row.Cells.Count; $c++) {\n    $row.Cells.Item($c).Range.Text = $values[$c - 1]\n  }\n}\n"}
